# Fix #215 - only report active conditions
#
# Applies, via the Word object model, the same textual change that the
# target unified diff makes to letter_to_request_inspection.docx:
#   1. Paragraph "I am writing to request ... state sanitary code. "
#      collapses into a single run (the stray _GoBack bookmark that used
#      to sit between "sanitary " and "code. " goes away from here).
#   2. The "{% for condition in ... .filter(category="emergency") %}"
#      line collapses its trailing two runs into one (text unchanged).
#   3. Both "{{ condition.description }} {{ condition.code }}" lines
#      (the "emergency" loop body and the "else" loop body) become
#      "{{ condition.description or condition.original_description }}
#      (C.M.R. {{ condition.code or "410.00" }})" -- note the first of
#      the two occurrences ends up with no space before the final "}}"
#      while the second keeps the space, matching the source diff
#      exactly. The _GoBack bookmark re-appears at the very end of the
#      second occurrence's paragraph.

function Set-RangeTextNoAutocorrect {
    # Assigns $NewText to the sub-string $OldText (first match) inside
    # paragraph $ParaIndex of $Document, without Word's Find/Replace
    # "smart quotes" autocorrect kicking in (straight double quotes in
    # $NewText must stay straight). We append a throwaway marker
    # character so Word is forced to actually rewrite/re-flow the run
    # rather than treat an unchanged portion as a no-op, then delete
    # the marker in a second pass.
    param($Document, $ParaIndex, $OldText, $NewText)

    $p = $Document.Paragraphs.Item($ParaIndex)
    $full = $p.Range.Text
    $idx = $full.IndexOf($OldText)
    if ($idx -lt 0) {
        throw "Set-RangeTextNoAutocorrect: could not find [$OldText] in paragraph $ParaIndex (actual text: [$full])"
    }
    $pstart = $p.Range.Start
    $rng = $Document.Range($pstart + $idx, $pstart + $idx + $OldText.Length)
    $rng.Text = $NewText + "@"

    $p2 = $Document.Paragraphs.Item($ParaIndex)
    $markerStart = $p2.Range.Start + $idx + $NewText.Length
    $markerRng = $Document.Range($markerStart, $markerStart + 1)
    $markerRng.Text = ""
}

$d = $word.ActiveDocument

# --- 1. Merge the "state sanitary " / bookmark / "code. " runs --------
$p1 = $d.Paragraphs.Item(23)
$rng1 = $p1.Range
$rng1.Find.Execute(
    "state sanitary code. ",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "state sanitary code. ", 1) | Out-Null

# --- 2. Merge the ".filter(...)" / " %}" runs (no text change) --------
Set-RangeTextNoAutocorrect $d 26 '.filter(category="emergency") %}' '.filter(category="emergency") %}'

# --- 3. First "{{ condition.description }} {{ condition.code }}" ------
#        (inside the emergency-only loop) -- no space before final }}
Set-RangeTextNoAutocorrect $d 27 `
    '{{ condition.description }} {{ condition.code }}' `
    '{{ condition.description or condition.original_description }} (C.M.R. {{ condition.code or "410.00"}})'

# --- 4. Second "{{ condition.description }} {{ condition.code }}" -----
#        (inside the else/all-conditions loop) -- space before final }}
Set-RangeTextNoAutocorrect $d 32 `
    '{{ condition.description }} {{ condition.code }}' `
    '{{ condition.description or condition.original_description }} (C.M.R. {{ condition.code or "410.00" }})'

# Move the _GoBack bookmark: delete it from wherever it still is now
# (it should already be gone -- it was removed by step 1 -- but be
# defensive) and re-add it, collapsed, at the very end of paragraph 32.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$p32 = $d.Paragraphs.Item(32)
$endRng = $p32.Range
$endRng.Collapse(0)
$endRng.MoveEnd(1, -1) | Out-Null
$endRng.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRng) | Out-Null

Write-Host "p23:" $d.Paragraphs.Item(23).Range.Text
Write-Host "p26:" $d.Paragraphs.Item(26).Range.Text
Write-Host "p27:" $d.Paragraphs.Item(27).Range.Text
Write-Host "p32:" $d.Paragraphs.Item(32).Range.Text
Write-Host "GoBack exists:" $d.Bookmarks.Exists("_GoBack")
